$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2353.8333
$ws.Range("I137").Value = 1305.1852
$ws.Range("J137").Value = 3702.0952
$ws.Range("K137").Value = 3915.5556
$ws.Range("L137").Value = 11106.2856
$ws.Range("M137").Value = -1365.5556
$ws.Range("N137").Value = -16206.2856
$ws.Range("H138").Value = 2122.66
$ws.Range("I138").Value = 674.6829
$ws.Range("J138").Value = 3128.8813
$ws.Range("K138").Value = 2024.0487
$ws.Range("L138").Value = 9386.643899999999
$ws.Range("M138").Value = 3115.9513
$ws.Range("N138").Value = -19666.6439
$ws.Range("H141").Value = 2824.4119
$ws.Range("I141").Value = 2608.3333
$ws.Range("J141").Value = 3944.0908
$ws.Range("K141").Value = 7824.999899999999
$ws.Range("L141").Value = 11832.2724
$ws.Range("M141").Value = -2644.999899999999
$ws.Range("N141").Value = -22192.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 433.5
$ws.Range("I5").Value = 420.2
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 420.2
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -308.2
$ws.Range("N5").Value = -724
$ws.Range("H32").Value = 3351.5784
$ws.Range("I32").Value = 3112.3835
$ws.Range("K32").Value = 3112.3835
$ws.Range("M32").Value = -2825.3835
$ws.Range("H61").Value = 909.89655
$ws.Range("I61").Value = 625.381
$ws.Range("J61").Value = 1656.75
$ws.Range("K61").Value = 625.381
$ws.Range("L61").Value = 1656.75
$ws.Range("M61").Value = -413.381
$ws.Range("N61").Value = -2080.75
$ws.Range("H132").Value = 2387.1914
$ws.Range("I132").Value = 1614.2424
$ws.Range("K132").Value = 4842.7272
$ws.Range("M132").Value = -2312.7272
$ws.Range("H136").Value = 909.89655
$ws.Range("I136").Value = 625.381
$ws.Range("J136").Value = 1656.75
$ws.Range("K136").Value = 1876.143
$ws.Range("L136").Value = 4970.25
$ws.Range("M136").Value = 673.857
$ws.Range("N136").Value = -10070.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 433.5
$ws.Range("I4").Value = 420.2
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 420.2
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -305.2
$ws.Range("N4").Value = -730
$ws.Range("H80").Value = 135.22223
$ws.Range("I80").Value = 83.72727
$ws.Range("J80").Value = 170.625
$ws.Range("K80").Value = 83.72727
$ws.Range("L80").Value = 170.625
$ws.Range("M80").Value = 914.27273
$ws.Range("N80").Value = -2166.625
$ws.Range("H83").Value = 135.22223
$ws.Range("I83").Value = 83.72727
$ws.Range("J83").Value = 170.625
$ws.Range("K83").Value = 418.63635
$ws.Range("L83").Value = 853.125
$ws.Range("M83").Value = 4573.36365
$ws.Range("N83").Value = -10837.125
$ws.Range("H107").Value = 1917.75
$ws.Range("I107").Value = 1946.2
$ws.Range("J107").Value = 1870.3334
$ws.Range("K107").Value = 1946.2
$ws.Range("L107").Value = 1870.3334
$ws.Range("M107").Value = -26.20000000000005
$ws.Range("N107").Value = -5710.3334
$ws.Range("H134").Value = 1720.7887
$ws.Range("I134").Value = 955.5862
$ws.Range("J134").Value = 5134.769
$ws.Range("K134").Value = 2866.7586
$ws.Range("L134").Value = 15404.307
$ws.Range("M134").Value = -331.7586000000001
$ws.Range("N134").Value = -20474.307
$ws.Range("H141").Value = 30000
$ws.Range("J141").Value = 30000
$ws.Range("L141").Value = 30000
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6495659.5
$ws.Range("I31").Value = 1101.1459
$ws.Range("J31").Value = 17245274
$ws.Range("K31").Value = 1101.1459
$ws.Range("L31").Value = 17245274
$ws.Range("M31").Value = -806.1459
$ws.Range("N31").Value = -17245864
$ws.Range("H34").Value = 6495659.5
$ws.Range("I34").Value = 1101.1459
$ws.Range("J34").Value = 17245274
$ws.Range("K34").Value = 1101.1459
$ws.Range("L34").Value = 17245274
$ws.Range("M34").Value = -899.1459
$ws.Range("N34").Value = -17245678
$ws.Range("H58").Value = 1098.09
$ws.Range("I58").Value = 1221.9078
$ws.Range("J58").Value = 706
$ws.Range("K58").Value = 1221.9078
$ws.Range("L58").Value = 706
$ws.Range("M58").Value = -1018.9078
$ws.Range("N58").Value = -1112
$ws.Range("H99").Value = 9529060
$ws.Range("I99").Value = 15388180
$ws.Range("K99").Value = 15388180
$ws.Range("M99").Value = -15386682
$ws.Range("H126").Value = 9529060
$ws.Range("I126").Value = 15388180
$ws.Range("K126").Value = 46164540
$ws.Range("M126").Value = -46162070
$ws.Range("H132").Value = 2999.75
$ws.Range("I132").Value = 2847.476
$ws.Range("J132").Value = 3456.5715
$ws.Range("K132").Value = 8542.428
$ws.Range("L132").Value = 10369.7145
$ws.Range("M132").Value = -6012.428
$ws.Range("N132").Value = -15429.7145
$ws.Range("H134").Value = 2703.7612
$ws.Range("I134").Value = 3321.3513
$ws.Range("J134").Value = 1942.0667
$ws.Range("K134").Value = 9964.053899999999
$ws.Range("L134").Value = 5826.2001
$ws.Range("M134").Value = -7429.053899999999
$ws.Range("N134").Value = -10896.2001
$ws.Range("H136").Value = 1098.09
$ws.Range("I136").Value = 1221.9078
$ws.Range("J136").Value = 706
$ws.Range("K136").Value = 3665.7234
$ws.Range("L136").Value = 2118
$ws.Range("M136").Value = -1115.7234
$ws.Range("N136").Value = -7218

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1298.721
$ws.Range("I5").Value = 311.9
$ws.Range("K5").Value = 935.6999999999999
$ws.Range("M5").Value = -823.6999999999999
$ws.Range("H68").Value = 12036.111
$ws.Range("I68").Value = 705
$ws.Range("K68").Value = 2115
$ws.Range("M68").Value = -1304
$ws.Range("H71").Value = 12036.111
$ws.Range("I71").Value = 705
$ws.Range("K71").Value = 6345
$ws.Range("M71").Value = -2289
$ws.Range("H113").Value = 550.37256
$ws.Range("I113").Value = 516.6429000000001
$ws.Range("J113").Value = 591.43475
$ws.Range("K113").Value = 1549.9287
$ws.Range("L113").Value = 1774.30425
$ws.Range("M113").Value = 620.0712999999998
$ws.Range("N113").Value = -6114.30425
$ws.Range("H122").Value = 3044.7354
$ws.Range("J122").Value = 3742.0386
$ws.Range("L122").Value = 33678.3474
$ws.Range("N122").Value = -38578.3474
$ws.Range("H135").Value = 1298.721
$ws.Range("I135").Value = 311.9
$ws.Range("K135").Value = 2807.1
$ws.Range("M135").Value = -272.0999999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2271.5122
$ws.Range("I132").Value = 1416.0667
$ws.Range("K132").Value = 4248.2001
$ws.Range("M132").Value = -1718.2001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3585595.2
$ws.Range("I93").Value = 7408373
$ws.Range("J93").Value = 1741.1875
$ws.Range("K93").Value = 7408373
$ws.Range("L93").Value = 1741.1875
$ws.Range("M93").Value = -7407125
$ws.Range("N93").Value = -4237.1875
$ws.Range("H122").Value = 7020.75
$ws.Range("I122").Value = 3311
$ws.Range("K122").Value = 9933
$ws.Range("M122").Value = -7483
$ws.Range("H136").Value = 2636.3171
$ws.Range("I136").Value = 1499.6451
$ws.Range("K136").Value = 4498.9353
$ws.Range("M136").Value = -1948.9353

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9260787
$ws.Range("I132").Value = 602.0714
$ws.Range("J132").Value = 41671436
$ws.Range("K132").Value = 1806.2142
$ws.Range("L132").Value = 125014308
$ws.Range("M132").Value = 723.7857999999999
$ws.Range("N132").Value = -125019368
$ws.Range("H141").Value = 42987.047
$ws.Range("J141").Value = 42987.047
$ws.Range("L141").Value = 42987.047
$ws.Range("N141").Value = -53347.047
